$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Apply the same style/format as the other header cells (e.g. H1) to the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF)
$dataI = @(2,5,1,7,1,1,8,6,8,8,7,8,4,2,8,1,6,4,4,4)
$dataJ = @(3,6,2,8,1,1,8,6,8,8,7,8,5,3,8,1,7,4,4,4)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
